$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 168.3077
$ws.Range("J53").Value = 250.16667
$ws.Range("L53").Value = 250.16667
$ws.Range("N53").Value = -1524.16667
$ws.Range("H92").Value = 632.76666
$ws.Range("J92").Value = 1963.125
$ws.Range("L92").Value = 1963.125
$ws.Range("N92").Value = -4459.125
$ws.Range("H96").Value = 277
$ws.Range("I96").Value = 385.2
$ws.Range("J96").Value = 186.83333
$ws.Range("K96").Value = 1155.6
$ws.Range("L96").Value = 560.49999
$ws.Range("M96").Value = 217.4000000000001
$ws.Range("N96").Value = -3306.49999
$ws.Range("H100").Value = 5087.579
$ws.Range("I100").Value = 2086.0908
$ws.Range("J100").Value = 9214.625
$ws.Range("K100").Value = 2086.0908
$ws.Range("L100").Value = 9214.625
$ws.Range("M100").Value = -1545.0908
$ws.Range("N100").Value = -10296.625
$ws.Range("H112").Value = 1351.1333
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1351.1333
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4053.3999
$ws.Range("M112").Value = ""
$ws.Range("N112").Value = -6269.3999
$ws.Range("H135").Value = 1591.5927
$ws.Range("I135").Value = 1611.4166
$ws.Range("K135").Value = 14502.7494
$ws.Range("M135").Value = -11967.7494

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2868.55
$ws.Range("I45").Value = 3429.4614
$ws.Range("J45").Value = 1826.8572
$ws.Range("K45").Value = 3429.4614
$ws.Range("L45").Value = 1826.8572
$ws.Range("M45").Value = -3052.4614
$ws.Range("N45").Value = -2580.8572
$ws.Range("H61").Value = 4043.423
$ws.Range("I61").Value = 1874.1111
$ws.Range("K61").Value = 1874.1111
$ws.Range("M61").Value = -1662.1111
$ws.Range("H74").Value = 4041.125
$ws.Range("I74").Value = 3607
$ws.Range("K74").Value = 3607
$ws.Range("M74").Value = -2733
$ws.Range("H77").Value = 4041.125
$ws.Range("I77").Value = 3607
$ws.Range("K77").Value = 18035
$ws.Range("M77").Value = -13667
$ws.Range("H97").Value = 690.43243
$ws.Range("I97").Value = 637.5185
$ws.Range("J97").Value = 833.3
$ws.Range("K97").Value = 637.5185
$ws.Range("L97").Value = 833.3
$ws.Range("M97").Value = -141.5185
$ws.Range("N97").Value = -1825.3
$ws.Range("H101").Value = 86546.71000000001
$ws.Range("J101").Value = 86546.71000000001
$ws.Range("L101").Value = 86546.71000000001
$ws.Range("N101").Value = -93036.71000000001
$ws.Range("H112").Value = 21045.666
$ws.Range("J112").Value = 21045.666
$ws.Range("L112").Value = 21045.666
$ws.Range("N112").Value = -23999.666
$ws.Range("H122").Value = 3059.8462
$ws.Range("I122").Value = 2763.7273
$ws.Range("J122").Value = 3443.0588
$ws.Range("K122").Value = 8291.1819
$ws.Range("L122").Value = 10329.1764
$ws.Range("M122").Value = -5841.1819
$ws.Range("N122").Value = -15229.1764
$ws.Range("H132").Value = 2491.0334
$ws.Range("I132").Value = 2388.6316
$ws.Range("K132").Value = 7165.8948
$ws.Range("M132").Value = -4635.8948
$ws.Range("H136").Value = 4043.423
$ws.Range("I136").Value = 1874.1111
$ws.Range("K136").Value = 5622.3333
$ws.Range("M136").Value = -3072.3333
$ws.Range("H139").Value = 79998.42999999999
$ws.Range("J139").Value = 79998.42999999999
$ws.Range("L139").Value = 79998.42999999999
$ws.Range("N139").Value = -90278.42999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = ""
$ws.Range("H54").Value = 8108.6665
$ws.Range("I54").Value = 7117.25
$ws.Range("K54").Value = 7117.25
$ws.Range("M54").Value = -6633.25
$ws.Range("H92").Value = 70398.164
$ws.Range("J92").Value = 70398.164
$ws.Range("L92").Value = 70398.164
$ws.Range("N92").Value = -75390.164
$ws.Range("H94").Value = 1037.3636
$ws.Range("I94").Value = 933.8214
$ws.Range("K94").Value = 933.8214
$ws.Range("M94").Value = -482.8214

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2141.3125
$ws.Range("I105").Value = 2543.111
$ws.Range("J105").Value = 1624.7142
$ws.Range("K105").Value = 2543.111
$ws.Range("L105").Value = 1624.7142
$ws.Range("M105").Value = -796.1109999999999
$ws.Range("N105").Value = -5118.7142
$ws.Range("H134").Value = 3248.5417
$ws.Range("I134").Value = 3153.2856
$ws.Range("J134").Value = 3915.3333
$ws.Range("K134").Value = 9459.856800000001
$ws.Range("L134").Value = 11745.9999
$ws.Range("M134").Value = -6924.856800000001
$ws.Range("N134").Value = -16815.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3140.889
$ws.Range("I5").Value = 966.8570999999999
$ws.Range("J5").Value = 10750
$ws.Range("K5").Value = 2900.5713
$ws.Range("L5").Value = 32250
$ws.Range("M5").Value = -2788.5713
$ws.Range("N5").Value = -32474
$ws.Range("H13").Value = 1349.5
$ws.Range("J13").Value = 2500
$ws.Range("L13").Value = 7500
$ws.Range("N13").Value = -7836
$ws.Range("H23").Value = 450.3
$ws.Range("I23").Value = 269.2
$ws.Range("K23").Value = 807.5999999999999
$ws.Range("M23").Value = -572.5999999999999
$ws.Range("H94").Value = 5840.8335
$ws.Range("I94").Value = 683.3333
$ws.Range("K94").Value = 2049.9999
$ws.Range("M94").Value = -1373.9999
$ws.Range("H102").Value = 6832.625
$ws.Range("I102").Value = 1850
$ws.Range("J102").Value = 9822.200000000001
$ws.Range("K102").Value = 5550
$ws.Range("L102").Value = 29466.6
$ws.Range("M102").Value = -3116
$ws.Range("N102").Value = -34334.60000000001
$ws.Range("H120").Value = 36026.4
$ws.Range("I120").Value = 20000
$ws.Range("K120").Value = 60000
$ws.Range("M120").Value = -55162
$ws.Range("H131").Value = 1840.9565
$ws.Range("I131").Value = 1028.75
$ws.Range("J131").Value = 2011.9474
$ws.Range("K131").Value = 3086.25
$ws.Range("L131").Value = 6035.8422
$ws.Range("M131").Value = 1953.75
$ws.Range("N131").Value = -16115.8422
$ws.Range("H135").Value = 3140.889
$ws.Range("I135").Value = 966.8570999999999
$ws.Range("J135").Value = 10750
$ws.Range("K135").Value = 8701.713899999999
$ws.Range("L135").Value = 96750
$ws.Range("M135").Value = -6166.713899999999
$ws.Range("N135").Value = -101820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2374.6956
$ws.Range("I102").Value = 1448.1177
$ws.Range("K102").Value = 1448.1177
$ws.Range("M102").Value = 173.8823
$ws.Range("H122").Value = 5820.074
$ws.Range("I122").Value = 3111.6
$ws.Range("J122").Value = 9205.666999999999
$ws.Range("K122").Value = 9334.799999999999
$ws.Range("L122").Value = 27617.001
$ws.Range("M122").Value = -6884.799999999999
$ws.Range("N122").Value = -32517.001
$ws.Range("H139").Value = 189979
$ws.Range("J139").Value = 189979
$ws.Range("L139").Value = 189979
$ws.Range("N139").Value = -200259

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5124.75
$ws.Range("I40").Value = 4780.375
$ws.Range("K40").Value = 4780.375
$ws.Range("M40").Value = -4644.375
$ws.Range("H69").Value = 41450.2
$ws.Range("J69").Value = 47581.5
$ws.Range("L69").Value = 47581.5
$ws.Range("N69").Value = -49203.5
$ws.Range("H72").Value = 41450.2
$ws.Range("J72").Value = 47581.5
$ws.Range("L72").Value = 142744.5
$ws.Range("N72").Value = -150856.5
$ws.Range("H93").Value = 2400.577
$ws.Range("I93").Value = 2316.3157
$ws.Range("J93").Value = 2629.2856
$ws.Range("K93").Value = 2316.3157
$ws.Range("L93").Value = 2629.2856
$ws.Range("M93").Value = -1068.3157
$ws.Range("N93").Value = -5125.2856
$ws.Range("H94").Value = 28216
$ws.Range("J94").Value = 28216
$ws.Range("L94").Value = 28216
$ws.Range("N94").Value = -29568
$ws.Range("H110").Value = 62081.125
$ws.Range("J110").Value = 62081.125
$ws.Range("L110").Value = 62081.125
$ws.Range("N110").Value = -70261.125
$ws.Range("H122").Value = 3164.8948
$ws.Range("I122").Value = 2855.1875
$ws.Range("J122").Value = 4816.6665
$ws.Range("K122").Value = 8565.5625
$ws.Range("L122").Value = 14449.9995
$ws.Range("M122").Value = -6115.5625
$ws.Range("N122").Value = -19349.9995
$ws.Range("H132").Value = 2762.7727
$ws.Range("I132").Value = 2643.6924
$ws.Range("J132").Value = 3691.6
$ws.Range("K132").Value = 7931.0772
$ws.Range("L132").Value = 11074.8
$ws.Range("M132").Value = -5401.0772
$ws.Range("N132").Value = -16134.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3043.158
$ws.Range("I122").Value = 2880
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 8640
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -6190
$ws.Range("N122").Value = -15400
$ws.Range("H132").Value = 3893.3428
$ws.Range("I132").Value = 2577.238
$ws.Range("J132").Value = 5867.5
$ws.Range("K132").Value = 7731.714
$ws.Range("L132").Value = 17602.5
$ws.Range("M132").Value = -5201.714
$ws.Range("N132").Value = -22662.5
